$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value, taken from the updated cryptos
# price/volume snapshot. Values are applied as literal text so that
# numeric-looking strings (e.g. "1.00", "236.66") keep their original
# formatting instead of being normalized by Excel's automatic type
# detection (which would turn "1.00" into 1).
$updates = [ordered]@{
    "D2" = "96.981.10"
    "E2" = "  +0.71%  "
    "D3" = "3.689.66"
    "E3" = "  +0.94%  "
    "E4" = "  +0.01%  "
    "D5" = "236.66"
    "E5" = "  -1.87%  "
    "D6" = "1.92"
    "E6" = "  +2.14%  "
    "D7" = "656.78"
    "E7" = "  -0.45%  "
    "E8" = "  +0.21%  "
    "E9" = "  -1.36%  "
    "E10" = "  +0.01%  "
    "D11" = "3.688.03"
    "E11" = "  +0.94%  "
    "D12" = "44.14"
    "E12" = "  -1.30%  "
    "E13" = "  +2.31%  "
    "D14" = "0.0000299"
    "E14" = "  +11.20%  "
    "D15" = "6.75"
    "E15" = "  +1.66%  "
    "D16" = "4.376.24"
    "E16" = "  +0.97%  "
    "D17" = "96.770.87"
    "E17" = "  +0.72%  "
    "D18" = "9.10"
    "E18" = "  +2.56%  "
    "D19" = "3.675.93"
    "E19" = "  -0.24%  "
    "D20" = "13.01"
    "E20" = "  +2.08%  "
    "D21" = "18.71"
    "E21" = "  +2.48%  "
    "E22" = "  -4.24%  "
    "D23" = "519.56"
    "E23" = "  -0.59%  "
    "E24" = "  -1.17%  "
    "E25" = "  +3.46%  "
    "D26" = "6.93"
    "E26" = "  +0.49%  "
    "E27" = "  +25.08%  "
    "D28" = "101.27"
    "E28" = "  -0.64%  "
    "D29" = "13.34"
    "E29" = "  +2.96%  "
    "D30" = "12.43"
    "E30" = "  +1.69%  "
    "D31" = "3.01"
    "E31" = "  -0.35%  "
    "D32" = "1.00"
    "E32" = "  -0.06%  "
    "E33" = "  +2.10%  "
    "D34" = "1.86"
    "E34" = "  +1.28%  "
    "D35" = "0.998"
    "E35" = "  -0.35%  "
    "D36" = "32.16"
    "E36" = "  -3.07%  "
    "D37" = "646.82"
    "E37" = "  +3.32%  "
    "E38" = "  +0.43%  "
    "D39" = "8.81"
    "E39" = "  +0.88%  "
    "E40" = "  +0.02%  "
    "D41" = "6.84"
    "E41" = "  +10.20%  "
    "D42" = "2.05"
    "E42" = "  +5.84%  "
    "B43" = "Kaspa"
    "C43" = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
    "D43" = "0.161"
    "E43" = "  +1.54%  "
    "B44" = "EnergySwap"
    "C44" = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
    "D44" = "40.44"
    "E44" = "  -11.03%  "
    "B45" = "Algorand"
    "C45" = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
    "D45" = "0.481"
    "E45" = "  +14.91%  "
    "B46" = "ARBITRUM"
    "C46" = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
    "D46" = "0.955"
    "E46" = "  -0.36%  "
    "D47" = "0.0463"
    "E47" = "  +2.50%  "
    "D48" = "23.63"
    "E48" = "  +0.11%  "
    "E49" = "  -0.26%  "
    "D50" = "8.61"
    "E50" = "  +1.34%  "
    "D51" = "3.52"
    "E51" = "  -1.38%  "
}

foreach ($ref in $updates.Keys) {
    $value = $updates[$ref]
    $cell = $ws.Range($ref)

    if ($value -match '^-?\d+(\.\d+)?$') {
        # Looks like a plain number - force text storage so the exact
        # string (including trailing/leading zeros) is preserved, then
        # restore the cell's original style so no formatting changes.
        $origStyle = $cell.Style
        $cell.NumberFormat = "@"
        $cell.Value = $value
        $cell.Style = $origStyle
    } else {
        $cell.Value = $value
    }
}
